# Generate Report for Handoff
#
# Six files (rows 7, 8, 9, 11, 13, 14 on the zh-cn / de-de report sheets,
# same rows on the Overview sheet) are being (re)handed off:
#   - their Priority is stamped "ht" (handoff type) on both the zh-cn and
#     de-de sheets;
#   - the shared "Latest Handoff Datetime" on the zh-cn sheet moves from
#     08:22:11 to 08:22:32;
#   - the shared "Latest Handoff Datetime" on the de-de sheet and the
#     "Latest HO Xliff Generate Date" on the Overview sheet move from
#     08:22:18 to 08:22:37 (these two happen to have shared the same text
#     before the edit, so they move together).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

# ---- Overview sheet: "Latest HO Xliff Generate Date" (col G) ----------
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-18 08:22:37"
}

# ---- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 8).Value = "2016-08-18 08:22:32"   # H: Latest Handoff Datetime
    $zhcn.Cells.Item($r, 5).Value = "ht"                    # E: Priority
}

# ---- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 8).Value = "2016-08-18 08:22:37"   # H: Latest Handoff Datetime
    $dede.Cells.Item($r, 5).Value = "ht"                    # E: Priority
}
